$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the DAYSTAGE values to reflect that BloodPressureReading refers
# to a DAYSTAGE rather than a prescriptionScheduleEntry
$ws.Range("B2").Value = "MORNING"
$ws.Range("B3").Value = "MIDDAY"
$ws.Range("B4").Value = "AFTERNOON"

# Update the active selection to match the new data range
$ws.Range("B2:B4").Select()

$wb.Save()
